$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SampleAudioWaveformPaths")

# Fix mislabeled "Emotion" values in column A so they match the actual
# audio/waveform file codes (ANG, DIS, SAD, FEA, NEU, HAP).
$ws.Range("A3").Value = "Disgust"
$ws.Range("A6").Value = "Neutral"
$ws.Range("A7").Value = "Happy"

$ws.Range("B10").Select()
